$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing values: remove "Planning and Production" from A3,
# shift SCM/Store/Factory Maintenance up one row, and append the new
# "Process" and "Production Planning and Control" department entries.
$ws.Range("A3").Value = "Process"
$ws.Range("A4").Value = "SCM"
$ws.Range("A5").Value = "Store"
$ws.Range("A6").Value = "Factory Maintenance"
$ws.Range("A7").Value = "Production Planning and Control"

# Select A3 to match the updated selection in the sheet view
$ws.Range("A3").Select()

# Widen column A to (best) fit the longer department names now stored
$ws.Columns.Item(1).ColumnWidth = 26.75
